# Updated symbol list on Sat Jan 14 21:48:04 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for the coin rows
# whose market data changed. Values are written as literal text (matching
# the workbook's existing inline-string cells) using a leading apostrophe
# so Excel doesn't coerce "303.65" / "4.75%" into numeric cells, then the
# cell style is reset to "Normal" so no stray number-format style sticks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Price (D), new Volume(1h) (E). $null means "leave D unchanged".
$updates = @(
    @{ Row = 2;  D = "303.65";    E = "4.75%" },
    @{ Row = 3;  D = "32.12";     E = "9.47%" },
    @{ Row = 4;  D = "5.257";     E = "-0.10%" },
    @{ Row = 5;  D = "0.07525";   E = "4.98%" },
    @{ Row = 6;  D = "7.895";     E = "5.82%" },
    @{ Row = 7;  D = $null;       E = "6.91%" },
    @{ Row = 8;  D = "1.497";     E = "7.17%" },
    @{ Row = 9;  D = "0.9216";    E = "1.24%" },
    @{ Row = 10; D = "0.1702";    E = "5.16%" },
    @{ Row = 11; D = "0.08010";   E = "4.59%" },
    @{ Row = 12; D = $null;       E = "3.72%" },
    @{ Row = 13; D = "0.03044";   E = "4.56%" },
    @{ Row = 14; D = "0.09910";   E = "9.92%" },
    @{ Row = 15; D = "0.001489";  E = "-6.27%" },
    @{ Row = 16; D = "0.04600";   E = "2.27%" },
    @{ Row = 17; D = "0.006462";  E = "5.47%" },
    @{ Row = 18; D = "3.461";     E = "-1.03%" },
    @{ Row = 19; D = "2.232";     E = "-0.01%" },
    @{ Row = 20; D = "0.3299";    E = "1.28%" },
    @{ Row = 21; D = $null;       E = "-1.18%" },
    @{ Row = 22; D = "4.474";     E = "11.54%" },
    @{ Row = 23; D = "0.1619";    E = "2.10%" },
    @{ Row = 24; D = "0.001217";  E = "1.13%" },
    @{ Row = 25; D = "0.004457";  E = "5.66%" },
    @{ Row = 26; D = "0.0001398"; E = "19.89%" },
    @{ Row = 27; D = $null;       E = "7.13%" },
    @{ Row = 39; D = "0.01705";   E = "2,530.71%" },
    @{ Row = 40; D = "0.04486";   E = "2.19%" },
    @{ Row = 41; D = "0.006961";  E = "-0.57%" },
    @{ Row = 42; D = "0.1351";    E = "6.43%" },
    @{ Row = 43; D = "0.002117";  E = "1.63%" },
    @{ Row = 44; D = "0.01281";   E = "-3.88%" },
    @{ Row = 45; D = "0.00006156"; E = "6.03%" },
    @{ Row = 46; D = "0.7089";    E = "-63.25%" },
    @{ Row = 47; D = "0.01498";   E = "16.34%" }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $ws.Cells.Item($row, 4).Value = "'" + $u.D
    }
    $ws.Cells.Item($row, 5).Value = "'" + $u.E

    # Reset formatting to default so the forced-text apostrophe entry
    # doesn't leave a quote-prefix / number-format style behind.
    $ws.Range($ws.Cells.Item($row, 4), $ws.Cells.Item($row, 5)).Style = "Normal"
}
